$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135; this shifts existing rows 135-196 down to 136-197
$ws.Rows.Item(135).EntireRow.Insert()

# Populate the newly inserted row 135 with the new data record
$ws.Range("A135").Value = 7
$ws.Range("B135").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C135").Value = "Ñuble"
$ws.Range("D135").Value = 44466
$ws.Range("D135").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E135").Value = 16
$ws.Range("F135").Value = 100114001
$ws.Range("G135").Value = "Papa"
$ws.Range("H135").Value = "Patagonia"
$ws.Range("I135").Value = "1a (guarda)"
$ws.Range("J135").Value = 160
$ws.Range("K135").Value = 9500
$ws.Range("L135").Value = 10000
$ws.Range("M135").Value = 9750
$ws.Range("N135").Value = "$/saco 25 kilos"
$ws.Range("O135").Value = "Región del Maule"
$ws.Range("P135").Value = 390
$ws.Range("Q135").Value = 25
$ws.Range("R135").Value = "Hortaliza"
